# Add Options processing for title and placeholder for other options
$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("Survey")
$options = $wb.Worksheets.Item("Options")

# Populate the Options sheet with a Title / value table
$options.Range("A1").Value = "Title"
$options.Range("A2").Value = "Number participants"
$options.Range("B1").Value = "Super survey"
$options.Range("B2").Value = 16

# Set column widths to match autosized "best fit" columns
# (closest attainable values to Excel's bestFit pixel-based widths of
# 19.29 / 12.43 characters for "Number participants" / "Super survey")
$options.Columns.Item(1).ColumnWidth = 18.451822916666668
$options.Columns.Item(2).ColumnWidth = 11.592447916666666

# Update selections to match the saved state
$survey.Range("H4").Select()
$options.Range("B2").Select()

# Make the Options sheet the active tab when the workbook is reopened
$options.Activate()

$wb.Save()
